$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (from H1, which carries the bold/border/
# centered header formatting) onto the two new header cells so I1/J1 end up
# sharing the same cell style as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows (plain, unstyled numeric cells like the existing H column)
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 3

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 5

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 4
